# Cotações atualizadas - 2025-10-04
# Append a new row of fund quotes to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 30

$ws.Cells.Item($row, 1).Value = 45934
$ws.Cells.Item($row, 2).Value = "21,4463"
$ws.Cells.Item($row, 3).Value = "15,2675"
$ws.Cells.Item($row, 4).Value = "15,4193"
$ws.Cells.Item($row, 5).Value = "15,4193"

$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
